$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the Job value in C2 (was "IT", now empty)
$ws.Range("C2").ClearContents()

# Fix the city in E2 (was "bucuresti", now "Ploiesti")
$ws.Range("E2").Value = "Ploiesti"

# Update the active selection to match the saved view state
$ws.Range("C9").Select() | Out-Null
